$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 206.8463543333333
$ws.Range("H2").Value = 620.5390629999999
$ws.Range("I2").Value = 0.9727792429746633
$ws.Range("J2").Value = 0.9727792429746634
$ws.Range("M2").Value = 47.78957466666667
$ws.Range("N2").Value = 143.368724
$ws.Range("O2").Value = 0.9275442943704633
$ws.Range("P2").Value = 0.927544294370463
$ws.Range("Q2").Value = 9885.099294940625
$ws.Range("R2").Value = 88965.89365446562
$ws.Range("S2").Value = 0.9022958365031676
$ws.Range("T2").Value = 0.9022958365031675
$ws.Range("G3").Value = 206.8463543333333
$ws.Range("H3").Value = 620.5390629999999
$ws.Range("I3").Value = 0.9727792429746633
$ws.Range("J3").Value = 0.9727792429746634
$ws.Range("M3").Value = 3.12343
$ws.Range("N3").Value = 9.370290000000001
$ws.Range("O3").Value = 0.06062242017370963
$ws.Range("P3").Value = 0.06062242017370962
$ws.Range("Q3").Value = 646.0701085153634
$ws.Range("R3").Value = 5814.63097663827
$ws.Range("S3").Value = 0.05897223200387321
$ws.Range("T3").Value = 0.05897223200387321
$ws.Range("G4").Value = 206.8463543333333
$ws.Range("H4").Value = 620.5390629999999
$ws.Range("I4").Value = 0.9727792429746633
$ws.Range("J4").Value = 0.9727792429746634
$ws.Range("O4").Value = 0.01183328545582722
$ws.Range("P4").Value = 0.01183328545582722
$ws.Range("Q4").Value = 126.1106369002249
$ws.Range("R4").Value = 1134.995732102024
$ws.Range("S4").Value = 0.0115111744676227
$ws.Range("T4").Value = 0.0115111744676227
$ws.Range("G5").Value = 3.181559666666666
$ws.Range("H5").Value = 9.544678999999999
$ws.Range("I5").Value = 0.01496258038481643
$ws.Range("J5").Value = 0.01496258038481643
$ws.Range("M5").Value = 47.78957466666667
$ws.Range("N5").Value = 143.368724
$ws.Range("O5").Value = 0.9275442943704633
$ws.Range("P5").Value = 0.927544294370463
$ws.Range("Q5").Value = 152.0453832466218
$ws.Range("R5").Value = 1368.408449219596
$ws.Range("S5").Value = 0.01387845606499589
$ws.Range("T5").Value = 0.01387845606499589
$ws.Range("G6").Value = 3.181559666666666
$ws.Range("H6").Value = 9.544678999999999
$ws.Range("I6").Value = 0.01496258038481643
$ws.Range("J6").Value = 0.01496258038481643
$ws.Range("M6").Value = 3.12343
$ws.Range("N6").Value = 9.370290000000001
$ws.Range("O6").Value = 0.06062242017370963
$ws.Range("P6").Value = 0.06062242017370962
$ws.Range("Q6").Value = 9.937378909656667
$ws.Range("S6").Value = 0.0009070678349712474
$ws.Range("T6").Value = 0.0009070678349712474
$ws.Range("G7").Value = 3.181559666666666
$ws.Range("H7").Value = 9.544678999999999
$ws.Range("I7").Value = 0.01496258038481643
$ws.Range("J7").Value = 0.01496258038481643
$ws.Range("O7").Value = 0.01183328545582722
$ws.Range("P7").Value = 0.01183328545582722
$ws.Range("Q7").Value = 1.939741781732444
$ws.Range("S7").Value = 0.0001770564848492939
$ws.Range("T7").Value = 0.0001770564848492939
$ws.Range("I8").Value = 0.01225817664052023
$ws.Range("J8").Value = 0.01225817664052023
$ws.Range("M8").Value = 47.78957466666667
$ws.Range("N8").Value = 143.368724
$ws.Range("O8").Value = 0.9275442943704633
$ws.Range("P8").Value = 0.927544294370463
$ws.Range("Q8").Value = 124.5640201942716
$ws.Range("R8").Value = 1121.076181748444
$ws.Range("S8").Value = 0.01137000180229983
$ws.Range("T8").Value = 0.01137000180229983
$ws.Range("I9").Value = 0.01225817664052023
$ws.Range("J9").Value = 0.01225817664052023
$ws.Range("M9").Value = 3.12343
$ws.Range("N9").Value = 9.370290000000001
$ws.Range("O9").Value = 0.06062242017370963
$ws.Range("P9").Value = 0.06062242017370962
$ws.Range("Q9").Value = 8.141252570443333
$ws.Range("R9").Value = 73.27127313399001
$ws.Range("S9").Value = 0.0007431203348651697
$ws.Range("T9").Value = 0.0007431203348651698
$ws.Range("I10").Value = 0.01225817664052023
$ws.Range("J10").Value = 0.01225817664052023
$ws.Range("O10").Value = 0.01183328545582722
$ws.Range("P10").Value = 0.01183328545582722
$ws.Range("S10").Value = 0.000145054503355229
$ws.Range("T10").Value = 0.000145054503355229
